$d = $word.ActiveDocument

# NOTE: We deliberately avoid passing the replacement text straight into
# Find.Execute's ReplaceWith parameter, because this runtime's "replace all"
# (wdReplaceAll) path applies smart-quote autocorrection (straight "'" becomes
# a curly "'"). Instead we locate the match with Find.Execute (no replacement)
# and then assign Range.Text directly, which leaves punctuation untouched.

function Replace-All($findText, $replaceText) {
    $rng = $d.Content
    $guard = 0
    while ($rng.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
        $rng.Text = $replaceText
        $rng.Collapse(0)
        $guard = $guard + 1
        if ($guard -gt 50) { break }
    }
}

function Replace-Occurrence($findText, $replaceText, $occurrence) {
    $rng = $d.Content
    $count = 0
    while ($rng.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
        $count = $count + 1
        if ($count -eq $occurrence) {
            $rng.Text = $replaceText
            return $true
        }
        $rng.Collapse(0)
    }
    return $false
}

# --- "English" -> "Anglais" (both occurrences: language-picker line & language heading) ---
Replace-All "English" "Anglais"

# --- language list line ---
Replace-All " / Portuguese / French / Thai / Vietnamese / Spanish" " / portugais / français / thaïlandais / vietnamien / espagnol"

# --- table: "Brief" / ":" -> "Résumé" / " :" ---
Replace-All "Brief" "Résumé"
$rng = $d.Content
$rng.Find.Execute("Résumé", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Collapse(0)
$rng.MoveEnd(1, 1) | Out-Null
$rng.Text = " :"

# --- table: "Target audience" / ":" -> "Public cible" / " :" ---
Replace-All "Target audience" "Public cible"
$rng = $d.Content
$rng.Find.Execute("Public cible", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Collapse(0)
$rng.MoveEnd(1, 1) | Out-Null
$rng.Text = " :"

# --- "We can't wait to meet you! " -> French ---
Replace-All "We can’t wait to meet you! " "Nous sommes impatients de vous rencontrer ! "

# --- "Hi " -> "Salut " ---
Replace-All "Hi " "Salut "

# --- "[PARTNER NAME]" -> "[NOM DU PARTENAIRE]" ---
Replace-All "[PARTNER NAME]" "[NOM DU PARTENAIRE]"

# --- "We hope you're as excited as we are for " -> French ---
Replace-All "We hope you’re as excited as we are for " "Nous espérons que vous êtes aussi enthousiastes que nous pour "

# --- Second "[EVENT NAME]" only (the first one, in the Subject line, stays in English) ---
Replace-Occurrence "[EVENT NAME]" "[NOM DE L'ÉVÉNEMENT]" 2 | Out-Null

# --- "In this email, we've linked/attached the following documents:" -> French ---
Replace-All "In this email, we’ve linked/attached the following documents:" "Nous avons joint les documents suivants au présent courriel :"

# --- "Your return flight tickets" -> French ---
Replace-All "Your return flight tickets" "Vos billets d'avion aller-retour ;"

# --- "Your accommodation booking details" -> French ---
Replace-All "Your accommodation booking details" "Les informations sur votre lieu d'hébergement ;"

# --- "If you have any questions, please contact us via " -> French ---
Replace-All "If you have any questions, please contact us via " "Si vous avez des questions, veuillez nous contacter par "

# --- "live chat" -> "chat en direct" ---
Replace-All "live chat" "chat en direct"

# --- First " or " (between live chat and WhatsApp links) -> " ou sur " ---
Replace-Occurrence " or " " ou sur " 1 | Out-Null

# --- "If you have any questions, please contact your country manager, " -> French ---
Replace-All "If you have any questions, please contact your country manager, " "Si vous avez des questions, veuillez contacter votre responsable national, "

# --- ", at " -> ", à l'adresse " ---
Replace-All ", at " ", à l'adresse "

# --- Remaining " or " (between EMAIL ADDRESS and WHATSAPP NO) -> " ou au" ---
Replace-Occurrence " or " " ou au" 1 | Out-Null

# --- "See you on the " -> "Rendez-vous le " ---
Replace-All "See you on the " "Rendez-vous le "

# --- "[DD]th" -> "[DD]" and the following "!" run -> " !" ---
Replace-All "[DD]th" "[DD]"
$rng = $d.Content
$rng.Find.Execute("[DD]", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Collapse(0)
$rng.MoveEnd(1, 1) | Out-Null
$rng.Text = " !"

# --- Comment text: "choose either one" -> "choisissez l'un ou l'autre" ---
for ($i = 1; $i -le $d.Comments.Count; $i++) {
    $c = $d.Comments.Item($i)
    if ($c.Range.Text -eq "choose either one") {
        $c.Range.Text = "choisissez l'un ou l'autre"
    }
}

Write-Output "done"
